# Loan RBI, Variable Instalments
#
# The "Repayment schedule" sheet gains a new (blank) column just before the
# existing "Late" column, shifting "Late", "heading"/"Paid Date" and
# "Outstanding" one column to the right (N -> O, O -> P, P -> Q). The new
# column inherits its width from the column immediately to its left. The
# sheet is then left as the active sheet/tab with the selection parked just
# past the used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("Repayment schedule")

# Insert a new blank column before column N ("Late"); this shifts the
# existing N/O/P columns (Late / Paid Date(heading) / Outstanding) right by
# one, to O/P/Q - matching how Excel's own "Insert Column" shifts data.
$ws.Columns("N").Insert()

# The freshly inserted column picks up the width of the column to its left
# (column M), just like Excel does when inserting a column.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab (it was "Transactions"
# before), and leave the selection just outside the populated range.
$ws.Activate()
[void]$ws.Range("R7").Select()
